# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 104
$ws1.Range("F4").Value = 120
$ws1.Range("F5").Value = 2779
$ws1.Range("F6").Value = 271
$ws1.Range("F7").Value = 389

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 104
$ws4.Range("F4").Value = 120
$ws4.Range("F5").Value = 2779
$ws4.Range("F6").Value = 271
$ws4.Range("F9").Value = 389
